# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns
# for rows 2-51. Price cells are written through a Text number format
# so purely numeric-looking strings (e.g. "246.34") stay text cells
# instead of being auto-coerced into numbers by Excel, matching the
# original inline-string cell type; ClearFormats() afterwards restores
# the cell to its original (default) style so no formatting changes leak in.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "42.611.03"
$cell.ClearFormats()
$ws.Range("E2").Value = "  +1.75%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.252.78"
$cell.ClearFormats()
$ws.Range("E3").Value = "  +0.96%  "
$ws.Range("E4").Value = "  +0.01%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "246.34"
$cell.ClearFormats()
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("E6").Value = "  +1.59%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "76.68"
$cell.ClearFormats()
$ws.Range("E7").Value = "  +1.47%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  +1.70%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "45.03"
$cell.ClearFormats()
$ws.Range("E10").Value = "  +11.50%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0957"
$cell.ClearFormats()
$ws.Range("E11").Value = "  +0.80%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "7.31"
$cell.ClearFormats()
$ws.Range("E12").Value = "  +2.31%  "
$ws.Range("E13").Value = "  -0.08%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "14.76"
$cell.ClearFormats()
$ws.Range("E14").Value = "  +0.08%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.866"
$cell.ClearFormats()
$ws.Range("E15").Value = "  +1.19%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "2.257.11"
$cell.ClearFormats()
$ws.Range("E16").Value = "  +0.95%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "42.433.44"
$cell.ClearFormats()
$ws.Range("E17").Value = "  +1.61%  "
$ws.Range("E18").Value = "  +4.56%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "6.22"
$cell.ClearFormats()
$ws.Range("E19").Value = "  +1.62%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "72.32"
$cell.ClearFormats()
$ws.Range("E20").Value = "  +1.10%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "11.19"
$cell.ClearFormats()
$ws.Range("E21").Value = "  +56.20%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "2.27"
$cell.ClearFormats()
$ws.Range("E22").Value = "  +0.02%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "232.81"
$cell.ClearFormats()
$ws.Range("E23").Value = "  +0.79%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "11.85"
$cell.ClearFormats()
$ws.Range("E24").Value = "  +3.37%  "
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("E26").Value = "  -1.43%  "
$ws.Range("E27").Value = "  +0.86%  "
$ws.Range("E28").Value = "  +3.66%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "167.66"
$cell.ClearFormats()
$ws.Range("E29").Value = "  -0.67%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "20.79"
$cell.ClearFormats()
$ws.Range("E30").Value = "  +1.38%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "5.89"
$cell.ClearFormats()
$ws.Range("E31").Value = "  +21.95%  "
$ws.Range("E32").Value = "  -0.87%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "31.36"
$cell.ClearFormats()
$ws.Range("E33").Value = "  -4.84%  "
$ws.Range("E34").Value = "  +0.86%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "4.80"
$cell.ClearFormats()
$ws.Range("E35").Value = "  +6.75%  "
$ws.Range("E36").Value = "  +0.65%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.0317"
$cell.ClearFormats()
$ws.Range("E37").Value = "  +5.98%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "14.21"
$cell.ClearFormats()
$ws.Range("E38").Value = "  +6.28%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "2.20"
$cell.ClearFormats()
$ws.Range("E39").Value = "  +1.61%  "
$ws.Range("E40").Value = "  -1.45%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "64.37"
$cell.ClearFormats()
$ws.Range("E41").Value = "  +6.55%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.204"
$cell.ClearFormats()
$ws.Range("E42").Value = "  +1.26%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "108.53"
$cell.ClearFormats()
$ws.Range("E43").Value = "  -2.75%  "
$ws.Range("E44").Value = "  +3.02%  "
$ws.Range("E45").Value = "  +2.18%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.998"
$cell.ClearFormats()
$ws.Range("E46").Value = "  +0.20%  "
$ws.Range("E47").Value = "  +8.81%  "
$ws.Range("E48").Value = "  +2.53%  "
$ws.Range("E49").Value = "  +2.92%  "
$ws.Range("E50").Value = "  +0.90%  "
$ws.Range("E51").Value = "  +1.06%  "
